$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataText")

# The APS (Annual Population Survey) derived rows (employment, self-employment,
# unemployment, inactivity - rates & volumes) all share the same "LatestPeriod"
# (column B) text. Update it from "Oct 2022-Sep 2023" to "Jan 2023-Dec 2023".
$newPeriod = "`t`nJan 2023-Dec 2023"

$ws.Range("B2").Value = $newPeriod
$ws.Range("B3").Value = $newPeriod
$ws.Range("B4").Value = $newPeriod
$ws.Range("B5").Value = $newPeriod
$ws.Range("B6").Value = $newPeriod
$ws.Range("B7").Value = $newPeriod
$ws.Range("B8").Value = $newPeriod
$ws.Range("B9").Value = $newPeriod

# Match the author's final selection/cursor position in the saved file.
$ws.Range("B9").Select() | Out-Null
